$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.904.31'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').Value = '3.441.87'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'584.66"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').Value = "'173.76"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.85%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '3.439.09'
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').Value = "'0.132"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.77%  '
$ws.Range('E12').Value = '  -1.81%  '
$ws.Range('D13').Value = '4.039.62'
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('E14').Value = '  +1.64%  '
$ws.Range('D15').Value = "'28.87"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.69%  '
$ws.Range('D16').Value = '65.882.28'
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('E17').Value = '  -0.65%  '
$ws.Range('D18').Value = '3.443.73'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('E19').Value = '  -1.12%  '
$ws.Range('D20').Value = "'13.79"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').Value = "'370.68"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.64%  '
$ws.Range('D22').Value = "'7.60"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.21%  '
$ws.Range('D23').Value = "'72.19"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.07%  '
$ws.Range('D24').Value = "'0.999"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = "'0.529"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('E26').Value = '  +3.24%  '
$ws.Range('D27').Value = "'9.71"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('E28').Value = '  +3.11%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = "'5.77"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.38%  '
$ws.Range('D31').Value = "'23.60"
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Value = "'1.98"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.79%  '
$ws.Range('D34').Value = "'1.28"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.81%  '
$ws.Range('D35').Value = "'6.99"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('E36').Value = '  +1.11%  '
$ws.Range('E37').Value = '  +1.48%  '
$ws.Range('D38').Value = "'0.878"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('D39').Value = "'28.39"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.64%  '
$ws.Range('D40').Value = "'1.79"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').Value = "'2.63"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('D42').Value = '2.771.84'
$ws.Range('E42').Value = '  +2.89%  '
$ws.Range('E43').Value = '  +0.51%  '
$ws.Range('D44').Value = "'6.44"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.26%  '
$ws.Range('D45').Value = "'0.0685"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = "'24.73"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.82%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = "'39.97"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.94%  '
$ws.Range('E48').Value = '  -1.13%  '
$ws.Range('D49').Value = "'324.00"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.08%  '
$ws.Range('E51').Value = '  +0.93%  '
